$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20-74 down to 21-75.
$ws.Rows.Item(20).EntireRow.Insert()

# Populate the new row 20 with the new daily price record.
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = "2021-09-14"
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 100112012
$ws.Range("G20").Value = "Espinaca"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 9000
$ws.Range("N20").Value = "$/docena de atados"
$ws.Range("O20").Value = "Región de La Araucanía"
$ws.Range("P20").Value = 3000
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Hortaliza"
